# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates column G ("K") for rows 2-58 with newly computed values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @(1,3,2,1,3,0,1,2,1,2,1,0,2,1,1,1,1,2,2,0,0,2,1,0,2,0,1,2,1,1,1,1,0,2,0,0,2,2,1,1,0,1,0,1,2,1,0,1,2,1,4,2,1,0,0,3,0)

$startRow = 2
for ($i = 0; $i -lt $newK.Length; $i++) {
    $row = $startRow + $i
    $ws.Range("G$row").Value = $newK[$i]
}
